# Updated cryptos list on Sat Sep 30 04:09:24 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table with
# the latest coinranking.com snapshot, and fixes the Stellar / BinanceUSD
# rows, which had swapped rank positions since the previous run.
#
# Note: several "Price" strings are plain decimals (e.g. "214.70", "1.00",
# "0.0890"). Excel's normal cell-input parser would silently coerce those to
# numbers (dropping trailing zeros / changing the cell type), so those
# particular writes are given a leading quote-prefix ('...) to force them to
# stay text -- exactly how the source sheet stores every Price/Volume cell.
# Prices that already contain two dots (e.g. "26.928.29") aren't valid
# numbers to Excel, so a plain assignment is enough for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '26.928.29'
$ws.Range('E2').Value = '  -0.09%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.671.15'
$ws.Range('E3').Value = '  +1.18%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.02%  '

# Row 5 - BNB
$ws.Range('D5').Value = "'214.70"
$ws.Range('E5').Value = '  -0.04%  '

# Row 6 - XRP
$ws.Range('D6').Value = "'0.518"
$ws.Range('E6').Value = '  +1.65%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  +0.00%  '

# Row 8 - Cardano
$ws.Range('E8').Value = '  +0.12%  '

# Row 9 - Dogecoin
$ws.Range('E9').Value = '  +0.58%  '

# Row 10 - Solana
$ws.Range('D10').Value = "'20.23"
$ws.Range('E10').Value = '  +0.00%  '

# Row 11 - TRON
$ws.Range('D11').Value = "'0.0890"
$ws.Range('E11').Value = '  +1.31%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.906.73'
$ws.Range('E12').Value = '  +1.17%  '

# Row 13 - WrappedEther
$ws.Range('D13').Value = '1.664.42'
$ws.Range('E13').Value = '  +0.74%  '

# Row 14 - Polkadot
$ws.Range('E14').Value = '  +0.10%  '

# Row 15 - Polygon
$ws.Range('D15').Value = "'0.525"
$ws.Range('E15').Value = '  +0.94%  '

# Row 16 - Litecoin
$ws.Range('D16').Value = "'65.46"
$ws.Range('E16').Value = '  +0.47%  '

# Row 17 - WrappedBTC
$ws.Range('D17').Value = '26.925.23'
$ws.Range('E17').Value = '  -0.12%  '

# Row 18 - Chainlink
$ws.Range('D18').Value = "'8.02"
$ws.Range('E18').Value = '  +3.47%  '

# Row 19 - BitcoinCash
$ws.Range('D19').Value = "'233.19"
$ws.Range('E19').Value = '  -1.10%  '

# Row 20 - ShibaInu
$ws.Range('E20').Value = '  -0.09%  '

# Row 21 - Dai
$ws.Range('E21').Value = '  +0.02%  '

# Row 22 - Uniswap
$ws.Range('D22').Value = "'4.41"
$ws.Range('E22').Value = '  -0.08%  '

# Row 23 - Avalanche
$ws.Range('E23').Value = '  -2.09%  '

# Row 24 - Toncoin
$ws.Range('E24').Value = '  -2.12%  '

# Row 25 - Monero
$ws.Range('D25').Value = "'145.65"
$ws.Range('E25').Value = '  +0.23%  '

# Row 26 - Cosmos
$ws.Range('E26').Value = '  -0.03%  '

# Row 27 - EthereumClassic
$ws.Range('D27').Value = "'15.90"
$ws.Range('E27').Value = '  +0.58%  '

# Row 28 / Row 29 - Stellar and BinanceUSD swapped rank position, plus
# refreshed price/volume figures.
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = "'0.112"
$ws.Range('E28').Value = '  -1.19%  '

$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.07%  '

# Row 30 - Hedera
$ws.Range('D30').Value = "'0.0498"
$ws.Range('E30').Value = '  +0.21%  '

# Row 31 - PancakeSwap
$ws.Range('E31').Value = '  -0.13%  '

# Row 32 - Filecoin
$ws.Range('E32').Value = '  +0.33%  '

# Row 33 - Maker
$ws.Range('D33').Value = '1.456.83'
$ws.Range('E33').Value = '  -6.32%  '

# Row 34 - InternetComputer(DFINITY)
$ws.Range('D34').Value = "'3.13"
$ws.Range('E34').Value = '  +1.37%  '

# Row 35 - LidoDAOToken
$ws.Range('E35').Value = '  +1.38%  '

# Row 36 - HuobiToken
$ws.Range('D36').Value = "'2.41"
$ws.Range('E36').Value = '  +0.04%  '

# Row 37 - ImmutableX
$ws.Range('D37').Value = "'0.579"
$ws.Range('E37').Value = '  -1.18%  '

# Row 38 - ARBITRUM
$ws.Range('D38').Value = "'0.898"
$ws.Range('E38').Value = '  +0.42%  '

# Row 39 - VeChain
$ws.Range('E39').Value = '  +0.81%  '

# Row 40 - WEMIXToken
$ws.Range('E40').Value = '  +13.50%  '

# Row 41 - FraxShare
$ws.Range('D41').Value = "'5.75"
$ws.Range('E41').Value = '  -4.08%  '

# Row 42 - PaxDollar
$ws.Range('E42').Value = '  +0.06%  '

# Row 43 - MXToken
$ws.Range('E43').Value = '  +2.97%  '

# Row 44 - Aave
$ws.Range('D44').Value = "'66.20"
$ws.Range('E44').Value = '  +0.18%  '

# Row 45 - RocketPoolETH
$ws.Range('D45').Value = '1.811.26'
$ws.Range('E45').Value = '  +1.08%  '

# Row 46 - TrustWalletToken
$ws.Range('E46').Value = '  +0.41%  '

# Row 47 - Quant
$ws.Range('D47').Value = "'90.41"
$ws.Range('E47').Value = '  +0.65%  '

# Row 48 - RenderToken
$ws.Range('E48').Value = '  +0.70%  '

# Row 49 - Algorand
$ws.Range('E49').Value = '  +2.60%  '

# Row 50 - Cronos
$ws.Range('E50').Value = '  +0.43%  '

# Row 51 - EnergySwap
$ws.Range('D51').Value = "'7.67"
$ws.Range('E51').Value = '  +0.62%  '
